# Applies the Contoso Supplier Agreement (Japanese) text revisions.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Northwind Traders は、Contoso のソフトドリンクとジュースの専用サプライヤーです。" `
             "Northwind Traders は、Contoso のソフト ドリンクとジュースの専属サプライヤーです。"

Replace-Text "一律料金" "遅延料金"

Replace-Text "1 か月あたり `$22.5" "1 か月あたり 1.5%"

Replace-Text "早期支払い割引" "早期支払割引"

Replace-Text "10 日以内の 2% 割引" "10 日以内 2% 割引"

Replace-Text "いずれかの当事者が有効期限の少なくとも 30 日前に終了の書面による通知を行わない限り、契約は別の年に自動的に更新されます。" `
             "いずれかの当事者が有効期限の少なくとも 30 日前に終了の書面による通知を行わない限り、契約はもう 1 年自動的に更新されます。"

Replace-Text "最小注文金額" "最小注文量"

Replace-Text "1 か月あたり 100 リリース" "1 か月あたり 100 ケース"

Replace-Text "最大注文金額" "最大注文量"

Replace-Text "20 時間/月" "1 か月あたり 500 ケース"
